$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row order (A1:W1)
$headers = @("ISO3","land","region","oecd","keine_daten","extraktivismus","fragilitaet","problematische_schuldenstruktur","vulnerabilitaet_naturkatastrophen","zahlungssituation","income","oeffentliche_schulden_bip","trend_oe_schulden_bip","oeffentliche_schulden_staatseinnahmen","trend_oe_schulden_staat","auslandsschulden_bip","trend_ausl_bip","auslandsschuldenstand_exporteinnahmen","trend_aus_schuldenstand_export","auslandsschuldendienst_exporteinnahmen","trend_ausl_schuldendienst_export","iwf_einschaetzung","link")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Country-level additions: column K (income group) and column W (link)
$rowData = @{
    3 = @("l", "https://erlassjahr.de/laenderinfos/afghanistan/")
    4 = @("lm", "https://erlassjahr.de/laenderinfos/angola/")
    5 = @("um", "https://erlassjahr.de/laenderinfos/albanien/")
    7 = @("um", "https://erlassjahr.de/laenderinfos/argentinien/")
    8 = @("um", "https://erlassjahr.de/laenderinfos/armenien/")
    9 = @("h", "https://erlassjahr.de/laenderinfos/antigua-barbuda/")
    13 = @("l", "https://erlassjahr.de/laenderinfos/burundi/")
    15 = @("l", "https://erlassjahr.de/laenderinfos/benin/")
    16 = @("l", "https://erlassjahr.de/laenderinfos/burkina-faso/")
    17 = @("lm", "https://erlassjahr.de/laenderinfos/bangladesch/")
    19 = @("h", "https://erlassjahr.de/laenderinfos/bahrain/")
    20 = @("h", "https://erlassjahr.de/laenderinfos/bahamas/")
    21 = @("um", "https://erlassjahr.de/laenderinfos/bosnien-herzigowina/")
    22 = @("um", "https://erlassjahr.de/laenderinfos/weissrussland/")
    23 = @("um", "https://erlassjahr.de/laenderinfos/belize/")
    24 = @("lm", "https://erlassjahr.de/laenderinfos/bolivien/")
    25 = @("um", "https://erlassjahr.de/laenderinfos/brasilien/")
    26 = @("h", "https://erlassjahr.de/laenderinfos/barbados/")
    28 = @("lm", "https://erlassjahr.de/laenderinfos/bhutan/")
    29 = @($null, "https://erlassjahr.de/laenderinfos/aserbaidschan/")
    30 = @("l", "https://erlassjahr.de/laenderinfos/zentralafrikanische-republik/")
    34 = @("um", "https://erlassjahr.de/laenderinfos/china/")
    35 = @("lm", "https://erlassjahr.de/laenderinfos/elfenbeinkueste/")
    36 = @("lm", "https://erlassjahr.de/laenderinfos/kamerun/")
    37 = @("l", "https://erlassjahr.de/laenderinfos/kongo-d-r/")
    38 = @("lm", "https://erlassjahr.de/laenderinfos/republik-kongo/")
    39 = @("um", "https://erlassjahr.de/laenderinfos/kolumbien/")
    40 = @("lm", "https://erlassjahr.de/laenderinfos/komoren/")
    41 = @("lm", "https://erlassjahr.de/laenderinfos/cabo-verde/")
    42 = @("um", "https://erlassjahr.de/laenderinfos/costa-rica/")
    47 = @("lm", "https://erlassjahr.de/laenderinfos/dschibuti/")
    48 = @("um", "https://erlassjahr.de/laenderinfos/dominica/")
    50 = @("um", "https://erlassjahr.de/laenderinfos/dominikanische-republik/")
    51 = @($null, "https://erlassjahr.de/laenderinfos/algerien/")
    52 = @("um", "https://erlassjahr.de/laenderinfos/ecuador/")
    53 = @("lm", "https://erlassjahr.de/laenderinfos/aegypten/")
    54 = @("l", "https://erlassjahr.de/laenderinfos/eritrea/")
    57 = @("l", "https://erlassjahr.de/laenderinfos/aethiopien/")
    59 = @($null, "https://erlassjahr.de/laenderinfos/fidschi/")
    61 = @("lm", "https://erlassjahr.de/laenderinfos/mikronesien/")
    62 = @("um", "https://erlassjahr.de/laenderinfos/gabun/")
    64 = @("um", "https://erlassjahr.de/laenderinfos/georgien/")
    65 = @("lm", "https://erlassjahr.de/laenderinfos/ghana/")
    66 = @("l", "https://erlassjahr.de/laenderinfos/guinea/")
    67 = @("l", "https://erlassjahr.de/laenderinfos/gambia/")
    68 = @("l", "https://erlassjahr.de/laenderinfos/guinea-bissau/")
    69 = @($null, "https://erlassjahr.de/laenderinfos/aequatorialguinea/")
    71 = @("um", "https://erlassjahr.de/laenderinfos/grenada/")
    72 = @("um", "https://erlassjahr.de/laenderinfos/guatemala/")
    73 = @("um", "https://erlassjahr.de/laenderinfos/guyana/")
    75 = @("lm", "https://erlassjahr.de/laenderinfos/honduras/")
    77 = @("l", "https://erlassjahr.de/laenderinfos/haiti/")
    79 = @("lm", "https://erlassjahr.de/laenderinfos/indonesien/")
    80 = @("lm", "https://erlassjahr.de/laenderinfos/indien/")
    82 = @("um", "https://erlassjahr.de/laenderinfos/iran/")
    83 = @($null, "https://erlassjahr.de/laenderinfos/irak/")
    87 = @("um", "https://erlassjahr.de/laenderinfos/jamaika/")
    88 = @("um", "https://erlassjahr.de/laenderinfos/jordanien/")
    90 = @("um", "https://erlassjahr.de/laenderinfos/kasachstan/")
    91 = @("lm", "https://erlassjahr.de/laenderinfos/kenia/")
    92 = @("lm", "https://erlassjahr.de/laenderinfos/kirgisistan/")
    93 = @("lm", "https://erlassjahr.de/laenderinfos/kambodscha/")
    94 = @("lm", "https://erlassjahr.de/laenderinfos/kiribati/")
    95 = @("h", "https://erlassjahr.de/laenderinfos/st.-kitts-nevis/")
    98 = @("lm", "https://erlassjahr.de/laenderinfos/laos/")
    99 = @("um", "https://erlassjahr.de/laenderinfos/libanon/")
    100 = @("l", "https://erlassjahr.de/laenderinfos/liberia/")
    102 = @("um", "https://erlassjahr.de/laenderinfos/st-lucia/")
    103 = @("um", "https://erlassjahr.de/laenderinfos/sri-lanka/")
    104 = @($null, "https://erlassjahr.de/laenderinfos/lesotho/")
    109 = @("lm", "https://erlassjahr.de/laenderinfos/marokko/")
    110 = @("lm", "https://erlassjahr.de/laenderinfos/moldau/")
    111 = @("l", "https://erlassjahr.de/laenderinfos/madagaskar/")
    112 = @("um", "https://erlassjahr.de/laenderinfos/malediven/")
    113 = @("um", "https://erlassjahr.de/laenderinfos/mexiko/")
    114 = @("um", "https://erlassjahr.de/laenderinfos/marshallinseln/")
    115 = @("um", "https://erlassjahr.de/laenderinfos/mazedonien/")
    116 = @("l", "https://erlassjahr.de/laenderinfos/mali/")
    118 = @("lm", "https://erlassjahr.de/laenderinfos/myanmar/")
    119 = @("um", "https://erlassjahr.de/laenderinfos/montenegro/")
    120 = @("lm", "https://erlassjahr.de/laenderinfos/mongolei/")
    121 = @("l", "https://erlassjahr.de/laenderinfos/mosambik/")
    122 = @("lm", "https://erlassjahr.de/laenderinfos/mauretanien/")
    123 = @("um", "https://erlassjahr.de/laenderinfos/mauritius/")
    124 = @("l", "https://erlassjahr.de/laenderinfos/malawi/")
    125 = @("um", "https://erlassjahr.de/laenderinfos/malaysia/")
    126 = @("um", "https://erlassjahr.de/laenderinfos/namibia/")
    127 = @("l", "https://erlassjahr.de/laenderinfos/niger/")
    128 = @("lm", "https://erlassjahr.de/laenderinfos/nigeria/")
    129 = @("lm", "https://erlassjahr.de/laenderinfos/nicaragua/")
    132 = @($null, "https://erlassjahr.de/laenderinfos/nepal/")
    133 = @("um", $null)
    135 = @("h", $null)
    136 = @("lm", "https://erlassjahr.de/laenderinfos/pakistan/")
    137 = @("h", "https://erlassjahr.de/laenderinfos/panama/")
    138 = @("um", "https://erlassjahr.de/laenderinfos/peru/")
    139 = @($null, "https://erlassjahr.de/laenderinfos/philippinen/")
    141 = @("lm", "https://erlassjahr.de/laenderinfos/papua-neuguinea/")
    146 = @("um", "https://erlassjahr.de/laenderinfos/paraguay/")
    148 = @($null, "https://erlassjahr.de/laenderinfos/kosovo/")
    151 = @("l", "https://erlassjahr.de/laenderinfos/ruanda/")
    153 = @("lm", "https://erlassjahr.de/laenderinfos/sudan/")
    154 = @("lm", "https://erlassjahr.de/laenderinfos/senegal/")
    155 = @("h", $null)
    156 = @($null, "https://erlassjahr.de/laenderinfos/salomonen/")
    157 = @("l", "https://erlassjahr.de/laenderinfos/sierra-leone/")
    158 = @("lm", "https://erlassjahr.de/laenderinfos/el-salvador/")
    160 = @($null, "https://erlassjahr.de/laenderinfos/somalia/")
    161 = @("um", "https://erlassjahr.de/laenderinfos/serbien/")
    162 = @("l", "https://erlassjahr.de/laenderinfos/suedsudan/")
    163 = @("lm", "https://erlassjahr.de/laenderinfos/sao-tome-principe/")
    164 = @("um", "https://erlassjahr.de/laenderinfos/surinam/")
    168 = @($null, "https://erlassjahr.de/laenderinfos/swasiland/")
    169 = @("h", "https://erlassjahr.de/laenderinfos/seychellen/")
    170 = @($null, "https://erlassjahr.de/laenderinfos/syrien/")
    171 = @("l", "https://erlassjahr.de/laenderinfos/tschad/")
    172 = @("l", "https://erlassjahr.de/laenderinfos/togo/")
    173 = @($null, "https://erlassjahr.de/laenderinfos/thailand/")
    174 = @("l", "https://erlassjahr.de/laenderinfos/tadschikistan/")
    175 = @("um", $null)
    176 = @($null, "https://erlassjahr.de/laenderinfos/ost-timor/")
    177 = @("um", "https://erlassjahr.de/laenderinfos/tonga/")
    178 = @($null, "https://erlassjahr.de/laenderinfos/trinidad-tobago/")
    179 = @("lm", "https://erlassjahr.de/laenderinfos/tunesien/")
    180 = @("um", "https://erlassjahr.de/laenderinfos/tuerkei/")
    181 = @("um", "https://erlassjahr.de/laenderinfos/tuvalu/")
    182 = @("h", "https://erlassjahr.de/laenderinfos/taiwan/")
    183 = @("l", "https://erlassjahr.de/laenderinfos/tansania/")
    184 = @("l", "https://erlassjahr.de/laenderinfos/uganda/")
    185 = @("lm", "https://erlassjahr.de/laenderinfos/ukraine/")
    186 = @("h", "https://erlassjahr.de/laenderinfos/uruguay/")
    188 = @($null, "https://erlassjahr.de/laenderinfos/usbekistan/")
    189 = @("um", "https://erlassjahr.de/laenderinfos/st-vincent/")
    190 = @("um", "https://erlassjahr.de/laenderinfos/venezuela/")
    191 = @("lm", "https://erlassjahr.de/laenderinfos/vietnam/")
    192 = @("lm", "https://erlassjahr.de/laenderinfos/vanuatu/")
    193 = @("um", "https://erlassjahr.de/laenderinfos/samoa/")
    194 = @("lm", "https://erlassjahr.de/laenderinfos/jemen/")
    195 = @("um", "https://erlassjahr.de/laenderinfos/suedafrika/")
    196 = @("lm", "https://erlassjahr.de/laenderinfos/sambia/")
    197 = @("lm", "https://erlassjahr.de/laenderinfos/simbabwe/")
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    if ($null -ne $vals[0]) { $ws.Cells.Item($r, 11).Value = $vals[0] }
    if ($null -ne $vals[1]) { $ws.Cells.Item($r, 23).Value = $vals[1] }
}

Write-Host "Done. Dimension should now be A1:W197"
